$d = $word.ActiveDocument

# 1. Fix typo "All reviwers" -> "All reviewers"
$d.Content.Find.Execute("All reviwers made", $true, $false, $false, $false, $false, $true, 1, $false, "All reviewers made", 2) | Out-Null

# 2. Update charity advisory-panel placeholders:
#    (NAME CHARITY) -> (!!!!!!NAME CHARITY!!!!!!!!)
#    (NOMINATED BY THE CHARITY) -> (nominated by the charity)
$d.Content.Find.Execute("(NAME CHARITY) and a survivor of domestic abuse (NOMINATED BY THE CHARITY).", $true, $false, $false, $false, $false, $true, 1, $false, "(!!!!!!NAME CHARITY!!!!!!!!) and a survivor of domestic abuse (nominated by the charity).", 2) | Out-Null

# 3. Replace the unspaced em dash with a spaced en dash around "again"
$d.Content.Find.Execute("associations" + [char]0x2014 + "again", $true, $false, $false, $false, $false, $true, 1, $false, "associations " + [char]0x2013 + " again", 2) | Out-Null

# 4. Insert "lived experiences of" before "victims in mind"
$d.Content.Find.Execute("done with the victims in mind", $true, $false, $false, $false, $false, $true, 1, $false, "done with the lived experiences of victims in mind", 2) | Out-Null

# 5. Remove the word "also" from "Reviewer 130837269 is also reluctant"
$d.Content.Find.Execute("Reviewer 130837269 is also reluctant", $true, $false, $false, $false, $false, $true, 1, $false, "Reviewer 130837269 is reluctant", 2) | Out-Null
